$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 172, shifting existing rows 172:190 down to 173:191.
$ws.Rows("172:172").Insert()

# Populate the newly inserted row with the new weekly price-point data.
$ws.Range("A172").Value = 7
$ws.Range("B172").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C172").Value = "Ñuble"
$ws.Range("D172").Value = 45180
$ws.Range("E172").Value = 16
$ws.Range("F172").Value = "Fruta"
$ws.Range("G172").Value = 100108
$ws.Range("H172").Value = "Tropicales y subtropicales"
$ws.Range("I172").Value = 100108002
$ws.Range("J172").Value = "Mango"
$ws.Range("K172").Value = "Sin especificar"
$ws.Range("L172").Value = "Primera"
$ws.Range("M172").Value = 50
$ws.Range("N172").Value = 13000
$ws.Range("O172").Value = 13000
$ws.Range("P172").Value = 13000
$ws.Range("Q172").Value = "$/bandeja 4 kilos"
$ws.Range("R172").Value = "Brasil"
$ws.Range("S172").Value = 3250
$ws.Range("T172").Value = 4
